# ApplianceData.xlsx edit: add "Event Record Status" lookup sheet + wire it
# into "Event Records" (new StatusText column/row) per the commit
# "Added Event Record Status values to ApplianceData spreadsheet".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new worksheet "Event Record Status" right before the
#    existing "Dispense Modification" sheet (so the final sheet order is
#    Data Overview, Shadow Data, Basic Ingest, Event Records,
#    Event Record Status, Dispense Modification).
# ---------------------------------------------------------------------
$dispenseSheet = $wb.Worksheets.Item("Dispense Modification")
$statusSheet = $wb.Worksheets.Add($dispenseSheet)
$statusSheet.Name = "Event Record Status"

# ---------------------------------------------------------------------
# 2) Populate the new sheet.
#    Row1 (merged B1:D1): title "Event Record Status"
#    Row2: headers "Value" | "Enumerated Value" | "Textual Interpretation (string)"
#    Rows3-34: hex code | firmware enum name | human readable text
# ---------------------------------------------------------------------

$statusRows = @(
    @("0x00","eNoError","Dispense Completed"),
    @("0x01","eUnknown_Error","Error: Unknown"),
    @("0x02","eTop_of_Tank_Error","Error: Top-of-Tank"),
    @("0x03","eCarbonator_Fill_Timeout_Error","Error: Carbonator Fill Timeout"),
    @("0x04","eOver_Pressure_Error","Error: Over Pressure"),
    @("0x05","eCarbonation_Timeout_Error","Error: Carbonation Timeout"),
    @("0x06","eError_Recovery_Brew","Error: Recovery Brew"),
    @("0x07","eHandle_Lift_Error","Error: Handle Lift"),
    @("0x08","ePuncture_Mechanism_Error","Error: Puncture Mechanism"),
    @("0x09","eCarbonation_Mechanism_Error","Error: Carbonation Mechanism"),
    @("0x80","eCleaning_Cycle_Completed","Cleaning Cycle Completed"),
    @("0x81","eRinsing_Cycle_Completed","Rinsing Cycle Completed"),
    @("0x82","eCO2_Module_Attached","CO2 Cylinder Attached"),
    @("0x83","eFirmware_Update_Passed","Firmware Update Passed"),
    @("0x84","eFirmware_Update_Failed","Firmware Update Failed"),
    @("0x85","eDrain_Cycle_Complete","Drain Cycle Completed"),
    @("0x86","eFreezeEventUpdate","Freeze Event Update"),
    @("0x87","eCritical_Error_OverTemp","Critical Error: OverTemp"),
    @("0x88","eCritical_Error_PuncMechFail","Critical Error: PuncMechFail"),
    @("0x89","eCritical_Error_TrickleFillTmout","Critical Error: TrickleFillTmout"),
    @("0x8A","eCritical_Error_ClnRinCWTFillTmout","Critical Error: ClnRinCWTFillTmout"),
    @("0x8B","eCritical_Error_ExtendedOPError","Critical Error: ExtendedOPError"),
    @("0x8C","eCritical_Error_BadMemClear","Critical Error: BadMemClear"),
    @("0xE0","eBLE_ModuleReset","BLE: ModuleReset"),
    @("0xE1","eBLE_IdleStatus","BLE: IdleStatus"),
    @("0xE2","eBLE_StandbyStatus","BLE: StandbyStatus"),
    @("0xE3","eBLE_ConnectedStatus","BLE: ConnectedStatus"),
    @("0xE4","eBLE_HealthTimeout","BLE: HealthTimeout"),
    @("0xE5","eBLE_ErrorState","BLE: ErrorState"),
    @("0xE6","eBLE_MultiConnectStat","BLE: MultiConnectStat"),
    @("0xE7","eBLE_MaxCriticalTimeout","BLE: MaxCriticalTimeout"),
    @("0xFF","eUnknownStatus","Unknown Status")
)

$statusSheet.Range("B1").Value = "Event Record Status"
$statusSheet.Range("B2").Value = "Value"
$statusSheet.Range("C2").Value = "Enumerated Value"
$statusSheet.Range("D2").Value = "Textual Interpretation (string)"

$r = 3
foreach ($row in $statusRows) {
    $statusSheet.Cells.Item($r, 2).Value = $row[0]
    $statusSheet.Cells.Item($r, 3).Value = $row[1]
    $statusSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

Write-Host "Populated Event Record Status rows"
